$d = $word.ActiveDocument

# Step 1: remove "Saavedra " after "Cervantes "
$d.Content.Find.Execute("Cervantes Saavedra fue", $true, $false, $false, $false, $false, $true, 1, $false, "Cervantes fue", 2) | Out-Null

# Step 2: replace "soldado español considerado" with "soldado nacido en España considerado"
$d.Content.Find.Execute("soldado español considerado", $true, $false, $false, $false, $false, $true, 1, $false, "soldado nacido en España considerado", 2) | Out-Null
